$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6, shifting existing rows 6-43 down to 7-44.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly record.
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 44685
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 100112040
$ws.Range("G6").Value = "Cilantro"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 550
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = 575
$ws.Range("N6").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O6").Value = "Provincia de Diguillín"
$ws.Range("P6").Value = 575
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = "Hortaliza"
